$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Vip -> Vipr1 -> ECs (self loop), now with updated numeric values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vip"
$ws.Range("C2").Value = "Vipr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.119963
$ws.Range("H2").Value = 3.359889
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1429513333333333
$ws.Range("N2").Value = 0.428854
$ws.Range("O2").Value = 0.1058099051556342
$ws.Range("P2").Value = 0.1058099051556342
$ws.Range("Q2").Value = 0.160100204134
$ws.Range("R2").Value = 1.440901837206
$ws.Range("S2").Value = 0.1058099051556342
$ws.Range("T2").Value = 0.1058099051556342

# Row 3: ECs -> Vip -> Vipr1 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vip"
$ws.Range("C3").Value = "Vipr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.119963
$ws.Range("H3").Value = 3.359889
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.027509
$ws.Range("N3").Value = 0.082527
$ws.Range("O3").Value = 0.02036164765346488
$ws.Range("P3").Value = 0.02036164765346487
$ws.Range("Q3").Value = 0.030809062167
$ws.Range("R3").Value = 0.277281559503
$ws.Range("S3").Value = 0.02036164765346488
$ws.Range("T3").Value = 0.02036164765346487

# Row 4: ECs -> Vip -> Vipr1 -> M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vip"
$ws.Range("C4").Value = "Vipr1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.119963
$ws.Range("H4").Value = 3.359889
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.047884
$ws.Range("N4").Value = 0.143652
$ws.Range("O4").Value = 0.03544284184225206
$ws.Range("P4").Value = 0.03544284184225206
$ws.Range("Q4").Value = 0.05362830829200001
$ws.Range("R4").Value = 0.482654774628
$ws.Range("S4").Value = 0.03544284184225206
$ws.Range("T4").Value = 0.03544284184225206

# Row 5: ECs -> Vip -> Vipr1 -> M2
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vip"
$ws.Range("C5").Value = "Vipr1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.119963
$ws.Range("H5").Value = 3.359889
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.9738796666666668
$ws.Range("N5").Value = 2.921639
$ws.Range("O5").Value = 0.7208475273379799
$ws.Range("P5").Value = 0.7208475273379799
$ws.Range("Q5").Value = 1.090709193119
$ws.Range("R5").Value = 9.816382738071001
$ws.Range("S5").Value = 0.7208475273379799
$ws.Range("T5").Value = 0.7208475273379799

# Row 6: ECs -> Vip -> Vipr1 -> sCs
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Vip"
$ws.Range("C6").Value = "Vipr1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.119963
$ws.Range("H6").Value = 3.359889
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1587963333333333
$ws.Range("N6").Value = 0.476389
$ws.Range("O6").Value = 0.1175380780106689
$ws.Range("P6").Value = 0.1175380780106689
$ws.Range("Q6").Value = 0.177846017869
$ws.Range("R6").Value = 1.600614160821
$ws.Range("S6").Value = 0.1175380780106689
$ws.Range("T6").Value = 0.1175380780106689
